$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper used to push literal percent-looking text ("40.9%") into a cell
# without Excel's normal autoconvert-to-number-with-percent-format
# behaviour clobbering the cell's existing (unrelated) style.
function Set-TextValue {
    param($address, $text)
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($address).PasteSpecial(-4163)
    $helper.Clear()
    $excel.CutCopyMode = 0
}

# --- "Recorded By" cleanup: drop stray "system/System, " prefix ---
$ws.Range("G2").Value  = "backup@backdoor.com"
$ws.Range("G6").Value  = "dnasr281@gmail.com"
$ws.Range("G12").Value = "dnasr281@gmail.com"
$ws.Range("G13").Value = "dnasr281@gmail.com"
$ws.Range("G29").Value = "backup@backdoor.com"
$ws.Range("G33").Value = "dnasr281@gmail.com"
$ws.Range("G39").Value = "dnasr281@gmail.com"
$ws.Range("G40").Value = "dnasr281@gmail.com"
$ws.Range("G56").Value = "backup@backdoor.com"
$ws.Range("G60").Value = "dnasr281@gmail.com"
$ws.Range("G66").Value = "dnasr281@gmail.com"
$ws.Range("G67").Value = "dnasr281@gmail.com"

# --- "Recorded By" reorder (admin/dnasr swapped) ---
$ws.Range("G90").Value  = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G116").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G142").Value = "dnasr281@gmail.com, admin@admin.com"

# --- Students (attendance count) updates ---
$ws.Range("H2").Value   = "32/53"
$ws.Range("H6").Value   = "43/53"
$ws.Range("H12").Value  = "30/53"
$ws.Range("H13").Value  = "35/53"
$ws.Range("H92").Value  = "45/56"
$ws.Range("H118").Value = "46/55"

# --- Class statistics (K/L columns) ---
$ws.Range("L6").Value  = 65
$ws.Range("L7").Value  = 4
Set-TextValue "L9"  "40.9%"
Set-TextValue "L10" "61.4%"

# --- Group statistics (M-S columns) ---
Set-TextValue "S15" "59.4%"

$ws.Range("O18").Value = 10
$ws.Range("P18").Value = 1
Set-TextValue "R18" "38.5%"
Set-TextValue "S18" "63.0%"

$ws.Range("O19").Value = 10
$ws.Range("P19").Value = 1
Set-TextValue "R19" "38.5%"
Set-TextValue "S19" "67.3%"

$ws.Range("O20").Value = 10
$ws.Range("P20").Value = 1
Set-TextValue "R20" "38.5%"
Set-TextValue "S20" "70.4%"

# --- Sessions flipped from Recorded to Not Recorded (copy the existing
#     "Not Recorded" look from row 31 so the pink style is reused, then
#     overwrite the text fields) ---
$ws.Range("A31:I31").Copy()

$ws.Range("A93:I93").PasteSpecial(-4122)
$ws.Range("G93").Value = ""
$ws.Range("H93").Value = "0/56"
$ws.Range("I93").Value = "Not Recorded"

$ws.Range("A119:I119").PasteSpecial(-4122)
$ws.Range("G119").Value = ""
$ws.Range("H119").Value = "0/55"
$ws.Range("I119").Value = "Not Recorded"

$ws.Range("A145:I145").PasteSpecial(-4122)
$ws.Range("G145").Value = ""
$ws.Range("H145").Value = "0/57"
$ws.Range("I145").Value = "Not Recorded"

$excel.CutCopyMode = 0
